# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" positioned between "总计" and "2022-Q2".
# - Populate it with the quarterly fund-holding table.
# - Update the "总计" (summary) sheet: the old 2022-Q2 summary row is pushed
#   down to row 3, and row 2 is rewritten with the new 2022-Q3 totals.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# ------------------------------------------------------------------
# 1. Push the existing 2022-Q2 summary row (row 2) down to row 3 on the
#    "总计" sheet, then overwrite row 2 with the new 2022-Q3 totals.
# ------------------------------------------------------------------
$oldB2 = $wsTotal.Range("B2").Value()
$oldC2 = $wsTotal.Range("C2").Value()
$oldD2 = $wsTotal.Range("D2").Value()

# Copy A2 (keeps its style) down to A3, then fix up its index value.
$wsTotal.Range("A2").Copy($wsTotal.Range("A3"))
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = $oldB2
$wsTotal.Range("C3").Value = $oldC2
$wsTotal.Range("D3").Value = $oldD2

# Row 2 becomes the new 2022-Q3 totals (A2 index/style stay as-is).
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0

# ------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right before "2022-Q2" so the
#    tab order becomes: 总计, 2022-Q3, 2022-Q2.
# ------------------------------------------------------------------
$wsQ2.Activate()
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q3"

# Match page margins used by the other data sheets in this workbook.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# ------------------------------------------------------------------
# 3. Header row - copy formatting (style 2) from the "总计" header and
#    then set the real header text.
# ------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $cols[$i] + "1"
    $wsTotal.Range("B1").Copy($newSheet.Range($cell))
    $newSheet.Range($cell).Value = $headers[$i]
}

# ------------------------------------------------------------------
# 4. Data rows. Column A is the numeric row-index column - copy its
#    style from the "总计" sheet. The numeric-looking text columns
#    (fund code / fund figures) are entered with a leading "'" so
#    Excel keeps them as literal text (e.g. "005167", "0.0031")
#    instead of parsing them as numbers, then the style is reset to
#    plain so no stray format is left on the cell.
# ------------------------------------------------------------------
$plainStyle = $wsTotal.Range("C2").Style

$wsTotal.Range("A2").Copy($newSheet.Range("A2"))
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'005167"
$newSheet.Range("C2").Value = "'嘉实润泽量化一年定期开放混合"
$newSheet.Range("D2").Value = "'0.55"
$newSheet.Range("E2").Value = "'24.55"
$newSheet.Range("F2").Value = "'0.57"
$newSheet.Range("G2").Value = "'0.0031"
$newSheet.Range("H2").Value = 4
$newSheet.Range("B2:G2").Style = $plainStyle

$wsTotal.Range("A2").Copy($newSheet.Range("A3"))
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'005166"
$newSheet.Range("C3").Value = "'嘉实润和量化6个月定期开放混合"
$newSheet.Range("D3").Value = "'0.22"
$newSheet.Range("E3").Value = "'24.64"
$newSheet.Range("F3").Value = "'0.55"
$newSheet.Range("G3").Value = "'0.0012"
$newSheet.Range("H3").Value = 4
$newSheet.Range("B3:G3").Style = $plainStyle

Write-Host "2022-Q3 sheet added"
